# Update market price data for Leve profit calculations across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 236.55556
$ws.Range("I33").Value = 203.93333
$ws.Range("J33").Value = 399.66666
$ws.Range("K33").Value = 203.93333
$ws.Range("L33").Value = 399.66666
$ws.Range("M33").Value = 25.06666999999999
$ws.Range("N33").Value = -857.66666
$ws.Range("H64").Value = 4701.143
$ws.Range("J64").Value = 4654.5
$ws.Range("L64").Value = 4654.5
$ws.Range("N64").Value = -5150.5
$ws.Range("H67").Value = 4701.143
$ws.Range("J67").Value = 4654.5
$ws.Range("L67").Value = 4654.5
$ws.Range("N67").Value = -6370.5
$ws.Range("H86").Value = 9551.166999999999
$ws.Range("J86").Value = 12351
$ws.Range("L86").Value = 12351
$ws.Range("N86").Value = -14597
$ws.Range("H88").Value = 3524.9
$ws.Range("I88").Value = 2541.7144
$ws.Range("J88").Value = 4054.3076
$ws.Range("K88").Value = 2541.7144
$ws.Range("L88").Value = 4054.3076
$ws.Range("M88").Value = -2135.7144
$ws.Range("N88").Value = -4866.3076
$ws.Range("H89").Value = 9551.166999999999
$ws.Range("J89").Value = 12351
$ws.Range("L89").Value = 61755
$ws.Range("N89").Value = -72987
$ws.Range("H91").Value = 3524.9
$ws.Range("I91").Value = 2541.7144
$ws.Range("J91").Value = 4054.3076
$ws.Range("K91").Value = 2541.7144
$ws.Range("L91").Value = 4054.3076
$ws.Range("M91").Value = -1137.7144
$ws.Range("N91").Value = -6862.3076
$ws.Range("H113").Value = 7999.25
$ws.Range("I113").Value = 7999.25
$ws.Range("K113").Value = 7999.25
$ws.Range("M113").Value = -4745.25
$ws.Range("H137").Value = 2168.886
$ws.Range("I137").Value = 1966.7797
$ws.Range("K137").Value = 5900.3391
$ws.Range("M137").Value = -3350.3391

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4909.4
$ws.Range("I63").Value = 4913.7144
$ws.Range("J63").Value = 4899.3335
$ws.Range("K63").Value = 4913.7144
$ws.Range("L63").Value = 4899.3335
$ws.Range("M63").Value = -4227.7144
$ws.Range("N63").Value = -6271.3335
$ws.Range("H66").Value = 4909.4
$ws.Range("I66").Value = 4913.7144
$ws.Range("J66").Value = 4899.3335
$ws.Range("K66").Value = 24568.572
$ws.Range("L66").Value = 24496.6675
$ws.Range("M66").Value = -21136.572
$ws.Range("N66").Value = -31360.6675
$ws.Range("H97").Value = 1351.3125
$ws.Range("I97").Value = 1351.3125
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1351.3125
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -855.3125
$ws.Range("N97").ClearContents()
$ws.Range("H110").Value = 2270
$ws.Range("I110").Value = 905
$ws.Range("K110").Value = 905
$ws.Range("M110").Value = 1140
$ws.Range("H122").Value = 3973.4443
$ws.Range("I122").Value = 3970.125
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 11910.375
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -9460.375
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3606.2354
$ws.Range("I94").Value = 3641
$ws.Range("J94").Value = 3493.25
$ws.Range("K94").Value = 3641
$ws.Range("L94").Value = 3493.25
$ws.Range("M94").Value = -3190
$ws.Range("N94").Value = -4395.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20002508
$ws.Range("I31").Value = 24392474
$ws.Range("J31").Value = 3774.889
$ws.Range("K31").Value = 24392474
$ws.Range("L31").Value = 3774.889
$ws.Range("M31").Value = -24392179
$ws.Range("N31").Value = -4364.889
$ws.Range("H34").Value = 20002508
$ws.Range("I34").Value = 24392474
$ws.Range("J34").Value = 3774.889
$ws.Range("K34").Value = 24392474
$ws.Range("L34").Value = 3774.889
$ws.Range("M34").Value = -24392272
$ws.Range("N34").Value = -4178.889
$ws.Range("H58").Value = 2425.3333
$ws.Range("J58").Value = 1999.6
$ws.Range("L58").Value = 1999.6
$ws.Range("N58").Value = -2405.6
$ws.Range("H107").Value = 1462.7028
$ws.Range("I107").Value = 1138.3793
$ws.Range("J107").Value = 2638.375
$ws.Range("K107").Value = 1138.3793
$ws.Range("L107").Value = 2638.375
$ws.Range("M107").Value = 781.6206999999999
$ws.Range("N107").Value = -6478.375
$ws.Range("H132").Value = 2541.5625
$ws.Range("I132").Value = 2510.6453
$ws.Range("K132").Value = 7531.9359
$ws.Range("M132").Value = -5001.9359
$ws.Range("H136").Value = 2425.3333
$ws.Range("J136").Value = 1999.6
$ws.Range("L136").Value = 5998.799999999999
$ws.Range("N136").Value = -11098.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6966
$ws.Range("I132").Value = 5615.6665
$ws.Range("J132").Value = 9666.666999999999
$ws.Range("K132").Value = 16846.9995
$ws.Range("L132").Value = 29000.001
$ws.Range("M132").Value = -14316.9995
$ws.Range("N132").Value = -34060.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5280.3335
$ws.Range("I7").Value = 5165.615
$ws.Range("J7").Value = 5578.6
$ws.Range("K7").Value = 5165.615
$ws.Range("L7").Value = 5578.6
$ws.Range("M7").Value = -5053.615
$ws.Range("N7").Value = -5802.6
$ws.Range("H22").Value = 56000
$ws.Range("J22").Value = 12000
$ws.Range("L22").Value = 12000
$ws.Range("N22").Value = -12590
$ws.Range("H27").Value = 56000
$ws.Range("J27").Value = 12000
$ws.Range("L27").Value = 12000
$ws.Range("N27").Value = -12214
$ws.Range("H46").Value = 3999.6667
$ws.Range("J46").Value = 3999.6667
$ws.Range("L46").Value = 3999.6667
$ws.Range("N46").Value = -4375.6667
$ws.Range("H126").Value = 5280.3335
$ws.Range("I126").Value = 5165.615
$ws.Range("J126").Value = 5578.6
$ws.Range("K126").Value = 15496.845
$ws.Range("L126").Value = 16735.8
$ws.Range("M126").Value = -13026.845
$ws.Range("N126").Value = -21675.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1307.8334
$ws.Range("I100").Value = 1393.4
$ws.Range("K100").Value = 2786.8
$ws.Range("M100").Value = -2245.8
$ws.Range("H132").Value = 1629.5
$ws.Range("I132").Value = 1600.7084
$ws.Range("J132").Value = 1975
$ws.Range("K132").Value = 4802.1252
$ws.Range("L132").Value = 5925
$ws.Range("M132").Value = -2272.1252
$ws.Range("N132").Value = -10985
